$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AWS")

# --- Re-point the "right aligned" cell style from D3 to D5/D9, leaving D3 plain ---

# D3 currently holds the only cell using style slot 2 (applyFill, cosmetically blank).
# Mutate alignment on D3 itself so the existing style slot turns into "horizontal=right".
$ws.Range("D3").HorizontalAlignment = -4152  # xlRight

# Propagate that style to D5 ("Nil") via a formats-only paste, so it reuses the
# same style slot instead of allocating a new one.
$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats

# New row 9: D9 = "18th " (right aligned, reusing the same style), E9 = "invocie- auto -dept"
$ws.Range("D9").Value = "18th "
$ws.Range("D3").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E9").Value = "invocie- auto -dept"

# Reset D3 back to the plain/default style (copy the unstyled format from D2).
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Other content changes ---

# Header B1: "Azure AD" -> "Azure AD Group"
$ws.Range("B1").Value = "Azure AD Group"

# Let column E size itself to fit the new content (bestFit-style autosize)
$ws.Columns.Item(5).ColumnWidth = 17.5

# Move the active selection to E12
$ws.Range("E12").Select()
